$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 77"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 78"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 79"
